$d = $word.ActiveDocument

# 1. Update the heading text
$d.Content.Find.Execute("Objet geolocalisation", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objet geoPositionUpdate", 2)

# 2. Update the table cell text "positionUpdate" -> "position"
$d.Content.Find.Execute("positionUpdate", $true, $false, $false, $false, $false,
                         $true, 1, $false, "position", 2)

# 3. Fill the empty "Exemple" cell (last cell, currently empty run) with a single space
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Count
$cell = $table.Cell($lastRow, 6)
$cell.Range.Text = " "
